$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 9 de Abril de 2020 a las 14:52"

# Update Galicia row (row 6) figures
$ws.Range("B6").Value = 6758
$ws.Range("C6").Value = 910
$ws.Range("D6").Value = 5551
$ws.Range("E6").Value = 297
